$d = $word.ActiveDocument
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("暫收支票收據列印(個人戶)")
Write-Host "Found:" $found
Write-Host "RangeStart:" $rng.Start "RangeEnd:" $rng.End
$rng.Text = "公平待客銀扣二扣資料"
